$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 7140
    $ws.Range("F4").Value = 5152
    $ws.Range("F11").Value = 86
    $ws.Range("F14").Value = 194
}
